$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.510.39"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.228.54"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.65"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "295.37"
$ws.Range("E6").Value = "  +10.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -1.39%  "

$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.88"
$ws.Range("E10").Value = "  -5.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.66"
$ws.Range("E13").Value = "  -6.37%  "

$ws.Range("E14").Value = "  +20.49%  "

$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.00"
$ws.Range("E16").Value = "  -2.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.566.06"
$ws.Range("E17").Value = "  -0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.233.50"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.559.94"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  +7.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("E21").Value = "  -1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.57"
$ws.Range("E22").Value = "  +2.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.37"
$ws.Range("E23").Value = "  +15.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.66"
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.90"
$ws.Range("E26").Value = "  -4.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.46"
$ws.Range("E28").Value = "  -6.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.27"
$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.53"
$ws.Range("E31").Value = "  -7.60%  "

$ws.Range("E32").Value = "  -4.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.40"
$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0884"
$ws.Range("E34").Value = "  -1.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  +1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.02"
$ws.Range("E36").Value = "  +7.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0375"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("E40").Value = "  -2.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.39"
$ws.Range("E41").Value = "  -4.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.20"
$ws.Range("E42").Value = "  +1.83%  "

$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.37"
$ws.Range("E45").Value = "  -6.99%  "

$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("E47").Value = "  -4.98%  "

$ws.Range("E48").Value = "  +3.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.66"
$ws.Range("E49").Value = "  +6.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.11"
$ws.Range("E50").Value = "  +2.07%  "

$ws.Range("E51").Value = "  +0.38%  "
